# Reorder the comma-separated "Recorded By" values in column G:
# reverse the order of the comma-separated entries in each cell that
# contains more than one entry (single-entry cells are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $text = $cell.Text

    if ($text -and $text.Contains(",")) {
        $parts = $text -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        $reversedParts = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $newText = [string]::Join(", ", $reversedParts)
        if (-not $newText.Equals($text)) {
            $cell.Value = $newText
        }
    }
}
